$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Value = 0.2305062539156956
$ws.Range("J8").Value = 0.2193215401759246
$ws.Range("I9").Value = 0.2109873117084238
$ws.Range("H10").Value = 0.2077622620068982
$ws.Range("G11").Value = 0.1868984584576193
$ws.Range("F12").Value = 0.2101374940836094
$ws.Range("E13").Value = 0.2201756597651073
$ws.Range("D14").Value = 0.1085991175498651
$ws.Range("C15").Value = 0.130019622424466
$ws.Range("B16").Value = 0.3662627537369125
